$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Prediction" (D) and "Error" (E) columns for rows 2-10, and
# D/E/F for row 11 (re-ran classification with updated model weights).
$ws.Range("D2").Value = 0.9999999243511539
$ws.Range("E2").Value = 0.9999999243511539

$ws.Range("D3").Value = 0.00000000003294034060409583
$ws.Range("E3").Value = 0.00000000003294034060409583

$ws.Range("D4").Value = 0.0003776924898923207
$ws.Range("E4").Value = 0.0003776924898923207

$ws.Range("D5").Value = 0.00000000002653290355792488
$ws.Range("E5").Value = 0.00000000002653290355792488

$ws.Range("D6").Value = 0.0000000000000000000000000000000000000000000002315719676356733
$ws.Range("E6").Value = 0.0000000000000000000000000000000000000000000002315719676356733

$ws.Range("D7").Value = 0.02668109272296164
$ws.Range("E7").Value = 0.9733189072770384

$ws.Range("D8").Value = 0.9999924035004871
$ws.Range("E8").Value = 0.000007596499512874111

$ws.Range("D9").Value = 0.9999582024683014
$ws.Range("E9").Value = 0.00004179753169863965

$ws.Range("D10").Value = 0.9999999227969958
$ws.Range("E10").Value = 0.00000007720300421176773

$ws.Range("D11").Value = 0.5223294696590324
$ws.Range("E11").Value = 0.4776705303409676
$ws.Range("F11").Value = 2.067084789276123
